$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: header + per-row "raw GOPS delta" (F - B), rounded to 2 dp.
$ws.Range("H1").Value = "Int8-SelfDependentW4A4"
for ($r = 2; $r -le 9; $r++) {
    $ws.Range("H$r").Formula = "=ROUND(F$r-B$r,2)"
}

# Match the new column's width to the rest of the wide label/data columns.
$ws.Columns.Item(8).ColumnWidth = 36.5

# Selection moved as part of the author's edit session.
$ws.Range("G20").Select() | Out-Null
